$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 6.265
$ws.Range("D3").Value = -7.345999999999999
$ws.Range("E8").Value = 16.703
$ws.Range("E11").Value = 17.107
$ws.Range("A12").Value = -21.544
$ws.Range("B14").Value = 6.132
$ws.Range("E14").Value = 16.991
$ws.Range("E15").Value = 16.185
$ws.Range("E17").Value = 16.725
$ws.Range("D20").Value = -7.765000000000001
$ws.Range("D25").Value = -8.196999999999999
$ws.Range("B26").Value = 5.920999999999999
$ws.Range("E26").Value = 16.966
$ws.Range("A27").Value = -21.531
$ws.Range("D30").Value = -7.37
$ws.Range("B31").Value = 6.425
$ws.Range("A32").Value = -21.345
$ws.Range("B35").Value = 7.287000000000001
$ws.Range("A36").Value = -21.343
$ws.Range("E36").Value = 16.64
$ws.Range("B37").Value = 7.848000000000002
$ws.Range("A38").Value = -20.093
$ws.Range("D44").Value = -8.161
$ws.Range("B45").Value = 5.928
$ws.Range("A46").Value = -21.481
$ws.Range("D47").Value = -7.575999999999999
$ws.Range("B52").Value = 5.378
$ws.Range("A54").Value = -21.862
$ws.Range("A55").Value = -22.21
$ws.Range("A56").Value = -22.097
$ws.Range("B57").Value = 5.332
$ws.Range("D58").Value = -8.178999999999998
$ws.Range("E64").Value = 17.185
$ws.Range("A67").Value = -21.603
$ws.Range("A69").Value = -21.636
$ws.Range("A72").Value = -21.567
$ws.Range("D78").Value = -7.812
$ws.Range("E79").Value = 17.266
$ws.Range("B81").Value = 6.375999999999999
$ws.Range("A83").Value = -20.146
$ws.Range("B83").Value = 7.326000000000001
$ws.Range("D84").Value = -8.196000000000002
$ws.Range("A86").Value = -22.172
$ws.Range("D89").Value = -7.105
$ws.Range("E89").Value = 17.517
$ws.Range("A91").Value = -21.526
$ws.Range("D91").Value = -7.139999999999999
$ws.Range("D92").Value = -7.056999999999999
$ws.Range("A93").Value = -21.665
$ws.Range("D96").Value = -7.654999999999999
$ws.Range("A99").Value = -20.437
$ws.Range("B100").Value = 5.517
$ws.Range("B102").Value = 7.499000000000001
$ws.Range("D102").Value = -7.865
